$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-10-04 18:29:34"
$ws.Range("A3").Value = "2025-10-04 18:29:34"
$ws.Range("A4").Value = "2025-10-04 18:29:34"
$ws.Range("A5").Value = "2025-10-04 18:29:34"
$ws.Range("A6").Value = "2025-10-04 18:29:34"
